$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header title (row 2): October 2016 -> November 2016
$ws.Range("A2").Value = "Commercial Sector by Census Division and State, Year-to-Date through November 2016"

# Row 4
$ws.Range("C4").Value = 129
$ws.Range("E4").Value = 30
$ws.Range("H4").Value = 662

# Row 5
$ws.Range("C5").Value = 587
$ws.Range("E5").Value = 60

# Row 6
$ws.Range("C6").Value = 549
$ws.Range("E6").Value = 221

# Row 7
$ws.Range("C7").Value = 79
$ws.Range("E7").Value = 28
$ws.Range("H7").Value = 662

# Row 8
$ws.Range("C8").Value = 252
$ws.Range("E8").Value = 241

# Row 9
$ws.Range("C9").Value = 233
$ws.Range("E9").Value = 175

# Row 10
$ws.Range("C10").Value = 1201
$ws.Range("E10").Value = 2815

# Row 11
$ws.Range("C11").Value = 107
$ws.Range("E11").Value = 28
$ws.Range("H11").Value = 775

# Row 12
$ws.Range("C12").Value = 486
$ws.Range("E12").Value = 82

# Row 13
$ws.Range("C13").Value = 98
$ws.Range("E13").Value = 29
$ws.Range("H13").Value = 775

# Row 14
$ws.Range("C14").Value = 1207
$ws.Range("E14").Value = 112

# Row 15
$ws.Range("B15").Value = 32
$ws.Range("C15").Value = 299
$ws.Range("E15").Value = 22
$ws.Range("H15").Value = 0

# Row 16
$ws.Range("B16").Value = 58
$ws.Range("C16").Value = 433
$ws.Range("E16").Value = 45
$ws.Range("H16").Value = 0

# Row 17
$ws.Range("B17").Value = 35
$ws.Range("C17").Value = 1039
$ws.Range("E17").Value = 46

# Row 18
$ws.Range("C18").Value = 241
$ws.Range("E18").Value = 28

# Row 19
$ws.Range("C19").Value = 1939
$ws.Range("E19").Value = 118

# Row 20
$ws.Range("C20").Value = 7275
$ws.Range("E20").Value = 93

# Row 21
$ws.Range("B21").Value = 36
$ws.Range("C21").Value = 113
$ws.Range("E21").Value = 48

# Row 22
$ws.Range("B22").Value = 75
$ws.Range("C22").Value = 607
$ws.Range("E22").Value = 108

# Row 23
$ws.Range("B23").Value = 423
$ws.Range("C23").Value = 119
$ws.Range("E23").Value = 87

# Row 24
$ws.Range("C24").Value = 510

# Row 25
$ws.Range("E25").Value = 1015

# Row 26
$ws.Range("C26").Value = 900

# Row 27
$ws.Range("C27").Value = 758

# Row 28
$ws.Range("B28").Value = 35
$ws.Range("C28").Value = 224
$ws.Range("E28").Value = 48
$ws.Range("H28").Value = 313

# Row 29
$ws.Range("C29").Value = 3091
$ws.Range("E29").Value = 156

# Row 30
$ws.Range("E30").Value = 148

# Row 31
$ws.Range("C31").Value = 53

# Row 32
$ws.Range("C32").Value = 1091
$ws.Range("E32").Value = 58

# Row 33
$ws.Range("C33").Value = 182
$ws.Range("H33").Value = 313

# Row 34
$ws.Range("C34").Value = 322
$ws.Range("E34").Value = 313
$ws.Range("H34").Value = 1081

# Row 35
$ws.Range("B35").Value = 244
$ws.Range("C35").Value = 181
$ws.Range("E35").Value = 442

# Row 36
$ws.Range("C36").Value = 810
$ws.Range("E36").Value = 95

# Row 37
$ws.Range("E37").Value = 348

# Row 38
$ws.Range("C38").Value = 810
$ws.Range("E38").Value = 99

# Row 39
$ws.Range("C39").Value = 1072
$ws.Range("E39").Value = 29

# Row 40
$ws.Range("E40").Value = 766

# Row 41
$ws.Range("E41").Value = 86

# Row 42
$ws.Range("C42").Value = 1512
$ws.Range("E42").Value = 253

# Row 43
$ws.Range("C43").Value = 1086
$ws.Range("E43").Value = 30

# Row 44
$ws.Range("C44").Value = 806
$ws.Range("E44").Value = 16
$ws.Range("H44").Value = 434

# Row 45
$ws.Range("C45").Value = 806
$ws.Range("E45").Value = 29

# Row 46
$ws.Range("H46").Value = 434

# Row 48
$ws.Range("E48").Value = 39

# Row 49
$ws.Range("E49").Value = 30

# Row 50
$ws.Range("E50").Value = 36

# Row 51
$ws.Range("C51").Value = 1409
$ws.Range("E51").Value = 9
$ws.Range("H51").Value = 330

# Row 52
$ws.Range("C52").Value = 1846
$ws.Range("E52").Value = 8
$ws.Range("H52").Value = 330

# Row 53
$ws.Range("C53").Value = 13658
$ws.Range("E53").Value = 128

# Row 54
$ws.Range("C54").Value = 211
$ws.Range("E54").Value = 218

# Row 55
$ws.Range("B55").Value = 37
$ws.Range("E55").Value = 1181

# Row 56
$ws.Range("B56").Value = 37
$ws.Range("C56").Value = 110
$ws.Range("E56").Value = 1181

# Row 58
$ws.Range("B58").Value = 20
$ws.Range("C58").Value = 77
$ws.Range("E58").Value = 9
$ws.Range("H58").Value = 201
